# This script refreshes the "IPC PO" model-prediction column (C) for the
# sliding-window results sheet after the DenseLayer/NeuralNetwork weight
# handling refactor produced new predictions, then recomputes the
# dependent DELTA (C-B), DELTA^2, TOTAL and MSE cells to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "IPC PO" (predicted) values for data rows 2-51, in row order.
$newIpcPo = @(30.08906540050407, 29.78556385453234, 29.2700097296191, 29.82163104762523, 29.50887915476766, 30.0574479233339, 29.93957089088588, 29.68718008544232, 29.61174624372592, 29.61864572433147, 29.69651013084202, 29.9117261406573, 30.20390420453156, 30.46567407497348, 30.22834259137082, 30.26397303977103, 30.78538835040724, 30.67079586766802, 30.81323862459335, 31.05625979610119, 31.21386468554301, 31.1276494698628, 31.05959569874924, 31.31717661401767, 31.71826208678449, 32.51659351030031, 32.51662862662089, 32.8520766034088, 33.22228952033154, 33.27917826661842, 33.34962080842794, 33.76405795937816, 33.76123245485297, 33.83241791371334, 34.22622927822735, 34.4503028803531, 35.51427435106902, 35.78944032712945, 36.07524671403826, 36.71440912003118, 36.80509697398391, 37.9356016899136, 38.71469226868126, 39.11247680791239, 39.5334708992003, 39.7926743465075, 40.15851467971533, 40.4026953359831, 40.73678562107821, 41.92965063556335)

$firstRow = 2
$lastRow = 51

$sumDelta = 0
$sumDeltaSq = 0

for ($i = 0; $i -lt $newIpcPo.Length; $i++) {
    $row = $firstRow + $i
    $ipcPo = $newIpcPo[$i]
    $ipcRo = $ws.Cells.Item($row, 2).Value2
    $delta = $ipcPo - $ipcRo
    $deltaSq = $delta * $delta

    $ws.Cells.Item($row, 3).Value2 = $ipcPo
    $ws.Cells.Item($row, 4).Value2 = $delta
    $ws.Cells.Item($row, 5).Value2 = $deltaSq

    $sumDelta += $delta
    $sumDeltaSq += $deltaSq
}

# TOTAL row: sum of DELTA in C52, sum of DELTA^2 in E52 (B52/D52 stay blank).
$ws.Cells.Item(52, 3).Value2 = $sumDelta
$ws.Cells.Item(52, 5).Value2 = $sumDeltaSq

# MSE row: mean of DELTA^2 in E53.
$rowCount = $lastRow - $firstRow + 1
$ws.Cells.Item(53, 5).Value2 = $sumDeltaSq / $rowCount
